# "add float content and anchor tags"
# Adds the new (3rd) attendance column E on the "3-4" sheet — rows 6-30 —
# with P(resent)/A(bsent) markers, mirroring the existing C/D columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-4")
$ws.Activate()

$attendance = @{
    6  = "A"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "P"
    13 = "P"
    14 = "P"
    15 = "A"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "A"
    24 = "P"
    25 = "A"
    26 = "P"
    27 = "P"
    28 = "A"
    29 = "P"
    30 = "P"
}

foreach ($row in 6..30) {
    $ws.Cells.Item($row, 5).Value = $attendance[$row]
}

# Anchor the view on the newly-entered data, same as Excel does after the
# last entry in a fill-down sequence.
$ws.Range("E31").Select()
